$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 14 ----
$ws.Range("A14").Value = 43692
$ws.Range("B14").Value = 0.5625
$ws.Range("C14").Value = "3 Hour"

$ws.Range("E14").Value = "Remove GDAL, divide large data file into smaller data files"

$ws.Range("F14").Value = "Update SOW with Weiwei"
$f14c = $ws.Range("F14").Characters(17, 6)
$f14c.Font.Bold = $true

$ws.Range("G14").Value = "Update SOW with Yuanxin"
$g14c = $ws.Range("G14").Characters(17, 7)
$g14c.Font.Bold = $true

$ws.Range("H14").Value = "Explore how to extract data from a WMS"

$ws.Range("I14").Value = "Trying to implement webserver with Peilin"
$i14c = $ws.Range("I14").Characters(36, 6)
$i14c.Font.Bold = $true

$ws.Range("J14").Value = "Trying to implement webserver with Dawei"
$j14c = $ws.Range("J14").Characters(36, 5)
$j14c.Font.Bold = $true

$ws.Range("K14").Value = "Research on maintainence cost"

# ---- Row 15 ----
$ws.Range("A15").Value = 43692
$ws.Range("B15").Value = 0.70833333333333337
$ws.Range("C15").Value = "1 Hour"

$ws.Range("H15").Value = "Explore how to extract data from a WMS. Work together with Yunyuan."
$h15c = $ws.Range("H15").Characters(60, 7)
$h15c.Font.Bold = $true

# E15's cell style changes from its original (border 17) to the interior
# border style (border 1) shared by F15/G15/H15/I15 - copy formats only
# from the untouched I15 cell (which already carries that style) before
# writing the new rich-text value into E15.
$ws.Range("I15").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E15").Value = "Explore how to extract data from a WMS. Work together with Daoyu."
$e15c = $ws.Range("E15").Characters(60, 5)
$e15c.Font.Bold = $true

# F15/G15 reuse the exact same rich strings as F14/G14 -> copy them across
# so the workbook shares the same sharedStrings entry (59/60) instead of
# minting new (non-deduped) rich-text entries.
$ws.Range("F14").Copy()
$ws.Range("F15").PasteSpecial(-4163)
$ws.Range("G14").Copy()
$ws.Range("G15").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# ---- Selection / view state ----
$av = $ws.Application.ActiveWindow
$av.ScrollRow = 10
$ws.Range("G15").Select()
